# Sync attendance_reports: rotate the "Recorded By" (column G) name/email
# list left by one entry for every data row where the column contains more
# than one comma-separated recorder (e.g. "System, dnasr281@gmail.com"
# becomes "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Column G = "Recorded By" (header in row 1); data starts at row 2.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($null -eq $current) {
        continue
    }

    $text = [string]$current
    if ($text.Contains(",")) {
        $parts = $text -split ",\s*"
        if ($parts.Length -gt 1) {
            $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
            $cell.Value = $rotated
        }
    }
}
